# Weekly CompStat data refresh: Manhattan North precinct report.
# Volume 31, Number 29 (week of 7/15/2024-7/21/2024)
#   -> Volume 31, Number 30 (week of 7/22/2024-7/28/2024)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: volume/issue number and date range ---
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Crime Complaints table (rows 14-33): refreshed weekly figures ---
# Row 14 (Murder)
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -50
$ws.Range("F14").Value = 7
$ws.Range("H14").Value = -12.5
$ws.Range("I14").Value = 31
$ws.Range("J14").Value = 36
$ws.Range("K14").Value = -13.888888888888
$ws.Range("L14").Value = 3.333333333333
$ws.Range("M14").Value = 3.333333333333
$ws.Range("N14").Value = -84.102564102564
# Row 15 (Rape)
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 9
$ws.Range("H15").Value = -10
$ws.Range("I15").Value = 82
$ws.Range("J15").Value = 76
$ws.Range("K15").Value = 7.894736842105
$ws.Range("L15").Value = -26.126126126126
$ws.Range("M15").Value = -28.695652173913
$ws.Range("N15").Value = -70.714285714285
# Row 16 (Robbery)
$ws.Range("C16").Value = 56
$ws.Range("D16").Value = 49
$ws.Range("E16").Value = 14.285714285714
$ws.Range("F16").Value = 161
$ws.Range("G16").Value = 164
$ws.Range("H16").Value = -1.829268292682
$ws.Range("I16").Value = 1198
$ws.Range("J16").Value = 1078
$ws.Range("K16").Value = 11.131725417439
$ws.Range("L16").Value = 1.870748299319
$ws.Range("M16").Value = -8.966565349544
$ws.Range("N16").Value = -77.794253938832
# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 44
$ws.Range("D17").Value = 75
$ws.Range("E17").Value = -41.333333333333
$ws.Range("F17").Value = 236
$ws.Range("G17").Value = 269
$ws.Range("H17").Value = -12.267657992565
$ws.Range("I17").Value = 1837
$ws.Range("J17").Value = 1740
$ws.Range("K17").Value = 5.574712643678
$ws.Range("L17").Value = 3.028603477285
$ws.Range("M17").Value = 66.094032549728
$ws.Range("N17").Value = -46.286549707602
# Row 18 (Burglary)
$ws.Range("C18").Value = 26
$ws.Range("D18").Value = 25
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 81
$ws.Range("G18").Value = 116
$ws.Range("H18").Value = -30.172413793103
$ws.Range("I18").Value = 767
$ws.Range("J18").Value = 897
$ws.Range("K18").Value = -14.492753623188
$ws.Range("L18").Value = -26.462128475551
$ws.Range("M18").Value = 0.130548302872
$ws.Range("N18").Value = -87.621045836023
# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 112
$ws.Range("D19").Value = 133
$ws.Range("E19").Value = -15.78947368421
$ws.Range("F19").Value = 473
$ws.Range("G19").Value = 523
$ws.Range("H19").Value = -9.560229445506
$ws.Range("I19").Value = 3565
$ws.Range("J19").Value = 3548
$ws.Range("K19").Value = 0.479143179255
$ws.Range("L19").Value = -3.33514099783
$ws.Range("M19").Value = 36.7472190257
$ws.Range("N19").Value = -42.509272697951
# Row 20 (G.L.A.)
$ws.Range("C20").Value = 20
$ws.Range("D20").Value = 24
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 74
$ws.Range("G20").Value = 112
$ws.Range("H20").Value = -33.928571428571
$ws.Range("I20").Value = 534
$ws.Range("J20").Value = 788
$ws.Range("K20").Value = -32.233502538071
$ws.Range("L20").Value = -28.418230563002
$ws.Range("M20").Value = 67.398119122257
$ws.Range("N20").Value = -90.0521609538
# Row 21 (TOTAL)
$ws.Range("C21").Value = 260
$ws.Range("D21").Value = 310
$ws.Range("E21").Value = -16.129032258064
$ws.Range("F21").Value = 1041
$ws.Range("G21").Value = 1202
$ws.Range("H21").Value = -13.394342762063
$ws.Range("I21").Value = 8014
$ws.Range("J21").Value = 8163
$ws.Range("K21").Value = -1.825309322552
$ws.Range("L21").Value = -6.564066689984
$ws.Range("M21").Value = 28.039622942962
$ws.Range("N21").Value = -70.378857882092
# Row 22 (Transit)
$ws.Range("C22").Value = 7
$ws.Range("E22").Value = 133.333333333333
$ws.Range("F22").Value = 16
$ws.Range("G22").Value = 20
$ws.Range("H22").Value = -20
$ws.Range("I22").Value = 139
$ws.Range("J22").Value = 167
$ws.Range("K22").Value = -16.766467065868
$ws.Range("L22").Value = -24.45652173913
$ws.Range("M22").Value = 5.30303030303
# Row 23 (Housing)
$ws.Range("C23").Value = 26
$ws.Range("D23").Value = 35
$ws.Range("E23").Value = -25.714285714285
$ws.Range("F23").Value = 100
$ws.Range("G23").Value = 112
$ws.Range("H23").Value = -10.714285714285
$ws.Range("I23").Value = 759
$ws.Range("J23").Value = 716
$ws.Range("K23").Value = 6.005586592178
$ws.Range("L23").Value = 2.016129032258
$ws.Range("M23").Value = 62.179487179487
# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 246
$ws.Range("D24").Value = 286
$ws.Range("E24").Value = -13.986013986014
$ws.Range("F24").Value = 1052
$ws.Range("G24").Value = 1145
$ws.Range("H24").Value = -8.122270742358
$ws.Range("I24").Value = 7293
$ws.Range("J24").Value = 8085
$ws.Range("K24").Value = -9.795918367346
$ws.Range("L24").Value = -19.351984960743
$ws.Range("M24").Value = 42.108339828526
# Row 25 (Retail Theft)
$ws.Range("C25").Value = 133
$ws.Range("D25").Value = 171
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 560
$ws.Range("G25").Value = 625
$ws.Range("H25").Value = -10.4
$ws.Range("I25").Value = 3876
$ws.Range("J25").Value = 4578
$ws.Range("K25").Value = -15.334207077326
$ws.Range("L25").Value = -32.497387669801
# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 94
$ws.Range("D26").Value = 90
$ws.Range("E26").Value = 4.444444444444
$ws.Range("F26").Value = 420
$ws.Range("G26").Value = 350
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 2962
$ws.Range("J26").Value = 2618
$ws.Range("K26").Value = 13.139801375095
$ws.Range("L26").Value = 12.537993920972
$ws.Range("M26").Value = -5.427841634738
# Row 27 (UCR Rape*)
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 15
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 130
$ws.Range("J27").Value = 149
$ws.Range("K27").Value = -12.751677852349
$ws.Range("L27").Value = -26.136363636363
# Row 28 (Other Sex Crimes)
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 13
$ws.Range("E28").Value = 23.076923076923
$ws.Range("F28").Value = 49
$ws.Range("G28").Value = 54
$ws.Range("H28").Value = -9.259259259259
$ws.Range("I28").Value = 333
$ws.Range("J28").Value = 345
$ws.Range("K28").Value = -3.478260869565
$ws.Range("L28").Value = -16.331658291457
# Row 29 (Shooting Vic.)
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 10
$ws.Range("E29").Value = -70
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = -66.666666666666
$ws.Range("I29").Value = 74
$ws.Range("J29").Value = 93
$ws.Range("K29").Value = -20.430107526881
$ws.Range("L29").Value = -37.81512605042
$ws.Range("M29").Value = -35.087719298245
$ws.Range("N29").Value = -84.518828451882
# Row 30 (Shooting Inc.)
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 9
$ws.Range("E30").Value = -66.666666666666
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 17
$ws.Range("H30").Value = -64.705882352941
$ws.Range("I30").Value = 60
$ws.Range("J30").Value = 84
$ws.Range("K30").Value = -28.571428571428
$ws.Range("L30").Value = -38.775510204081
$ws.Range("M30").Value = -39.393939393939
$ws.Range("N30").Value = -86.301369863013
# Row 31 (Hate Crimes)
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = 11
$ws.Range("H31").Value = -45.454545454545
$ws.Range("I31").Value = 66
$ws.Range("J31").Value = 48
$ws.Range("K31").Value = 37.5
$ws.Range("L31").Value = 4.761904761904
# Row 33 (Traffic Fatalities)
$ws.Range("F33").Value = 3
$ws.Range("H33").Value = 50
$ws.Range("I33").Value = 13
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = -7.142857142857
